$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.034777666666667
$ws.Range("H2").Value = 9.104333
$ws.Range("I2").Value = 0.2502264227183869
$ws.Range("J2").Value = 0.2502264227183869
$ws.Range("M2").Value = 15.50220733333333
$ws.Range("N2").Value = 46.506622
$ws.Range("O2").Value = 0.5994675913188158
$ws.Range("P2").Value = 0.5994675913188158
$ws.Range("Q2").Value = 47.04575259923622
$ws.Range("R2").Value = 423.411773393126
$ws.Range("S2").Value = 0.1500026309113152
$ws.Range("T2").Value = 0.1500026309113152

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.034777666666667
$ws.Range("H3").Value = 9.104333
$ws.Range("I3").Value = 0.2502264227183869
$ws.Range("J3").Value = 0.2502264227183869
$ws.Range("O3").Value = 0.04399860030713892
$ws.Range("P3").Value = 0.04399860030713892
$ws.Range("Q3").Value = 3.452976098688667
$ws.Range("R3").Value = 31.076784888198
$ws.Range("S3").Value = 0.01100961235947149
$ws.Range("T3").Value = 0.01100961235947149

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.034777666666667
$ws.Range("H4").Value = 9.104333
$ws.Range("I4").Value = 0.2502264227183869
$ws.Range("J4").Value = 0.2502264227183869
$ws.Range("M4").Value = 8.848210666666667
$ws.Range("N4").Value = 26.544632
$ws.Range("O4").Value = 0.3421587275782868
$ws.Range("P4").Value = 0.3421587275782868
$ws.Range("Q4").Value = 26.85235212116178
$ws.Range("R4").Value = 241.671169090456
$ws.Range("S4").Value = 0.08561715440378978
$ws.Range("T4").Value = 0.08561715440378978

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.034777666666667
$ws.Range("H5").Value = 9.104333
$ws.Range("I5").Value = 0.2502264227183869
$ws.Range("J5").Value = 0.2502264227183869
$ws.Range("M5").Value = 0.371739
$ws.Range("N5").Value = 1.115217
$ws.Range("O5").Value = 0.01437508079575842
$ws.Range("P5").Value = 0.01437508079575841
$ws.Range("Q5").Value = 1.128145215029
$ws.Range("R5").Value = 10.153306935261
$ws.Range("S5").Value = 0.003597025043810411
$ws.Range("T5").Value = 0.00359702504381041

# Row 6
$ws.Range("I6").Value = 0.4835045831069426
$ws.Range("J6").Value = 0.4835045831069426
$ws.Range("M6").Value = 15.50220733333333
$ws.Range("N6").Value = 46.506622
$ws.Range("O6").Value = 0.5994675913188158
$ws.Range("P6").Value = 0.5994675913188158
$ws.Range("Q6").Value = 90.90501614630088
$ws.Range("R6").Value = 818.145145316708
$ws.Range("S6").Value = 0.2898453278267271
$ws.Range("T6").Value = 0.2898453278267271

# Row 7
$ws.Range("I7").Value = 0.4835045831069426
$ws.Range("J7").Value = 0.4835045831069426
$ws.Range("O7").Value = 0.04399860030713892
$ws.Range("P7").Value = 0.04399860030713892
$ws.Range("S7").Value = 0.0212735248987922
$ws.Range("T7").Value = 0.0212735248987922

# Row 8
$ws.Range("I8").Value = 0.4835045831069426
$ws.Range("J8").Value = 0.4835045831069426
$ws.Range("M8").Value = 8.848210666666667
$ws.Range("N8").Value = 26.544632
$ws.Range("O8").Value = 0.3421587275782868
$ws.Range("P8").Value = 0.3421587275782868
$ws.Range("Q8").Value = 51.88594864098311
$ws.Range("R8").Value = 466.973537768848
$ws.Range("S8").Value = 0.1654353129341415
$ws.Range("T8").Value = 0.1654353129341415

# Row 9
$ws.Range("I9").Value = 0.4835045831069426
$ws.Range("J9").Value = 0.4835045831069426
$ws.Range("M9").Value = 0.371739
$ws.Range("N9").Value = 1.115217
$ws.Range("O9").Value = 0.01437508079575842
$ws.Range("P9").Value = 0.01437508079575841
$ws.Range("Q9").Value = 2.179879230782
$ws.Range("R9").Value = 19.618913077038
$ws.Range("S9").Value = 0.00695041744728179
$ws.Range("T9").Value = 0.006950417447281788

# Row 10
$ws.Range("G10").Value = 2.564975
$ws.Range("H10").Value = 7.694925
$ws.Range("I10").Value = 0.2114897989601526
$ws.Range("J10").Value = 0.2114897989601526
$ws.Range("M10").Value = 15.50220733333333
$ws.Range("N10").Value = 46.506622
$ws.Range("O10").Value = 0.5994675913188158
$ws.Range("P10").Value = 0.5994675913188158
$ws.Range("Q10").Value = 39.76277425481666
$ws.Range("R10").Value = 357.86496829335
$ws.Range("S10").Value = 0.1267812803711433
$ws.Range("T10").Value = 0.1267812803711433

# Row 11
$ws.Range("G11").Value = 2.564975
$ws.Range("H11").Value = 7.694925
$ws.Range("I11").Value = 0.2114897989601526
$ws.Range("J11").Value = 0.2114897989601526
$ws.Range("O11").Value = 0.04399860030713892
$ws.Range("P11").Value = 0.04399860030713892
$ws.Range("Q11").Value = 2.91843368495
$ws.Range("R11").Value = 26.26590316455
$ws.Range("S11").Value = 0.009305255133484921
$ws.Range("T11").Value = 0.009305255133484921

# Row 12
$ws.Range("G12").Value = 2.564975
$ws.Range("H12").Value = 7.694925
$ws.Range("I12").Value = 0.2114897989601526
$ws.Range("J12").Value = 0.2114897989601526
$ws.Range("M12").Value = 8.848210666666667
$ws.Range("N12").Value = 26.544632
$ws.Range("O12").Value = 0.3421587275782868
$ws.Range("P12").Value = 0.3421587275782868
$ws.Range("Q12").Value = 22.69543915473333
$ws.Range("R12").Value = 204.2589523926
$ws.Range("S12").Value = 0.0723630805079935
$ws.Range("T12").Value = 0.0723630805079935

# Row 13
$ws.Range("G13").Value = 2.564975
$ws.Range("H13").Value = 7.694925
$ws.Range("I13").Value = 0.2114897989601526
$ws.Range("J13").Value = 0.2114897989601526
$ws.Range("M13").Value = 0.371739
$ws.Range("N13").Value = 1.115217
$ws.Range("O13").Value = 0.01437508079575842
$ws.Range("P13").Value = 0.01437508079575841
$ws.Range("Q13").Value = 0.9535012415249999
$ws.Range("R13").Value = 8.581511173725
$ws.Range("S13").Value = 0.003040182947530898
$ws.Range("T13").Value = 0.003040182947530898

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.6643690000000001
$ws.Range("H14").Value = 1.993107
$ws.Range("I14").Value = 0.05477919521451775
$ws.Range("J14").Value = 0.05477919521451775
$ws.Range("M14").Value = 15.50220733333333
$ws.Range("N14").Value = 46.506622
$ws.Range("O14").Value = 0.5994675913188158
$ws.Range("P14").Value = 0.5994675913188158
$ws.Range("Q14").Value = 10.29918598383933
$ws.Range("R14").Value = 92.692673854554
$ws.Range("S14").Value = 0.03283835220963015
$ws.Range("T14").Value = 0.03283835220963015

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.6643690000000001
$ws.Range("H15").Value = 1.993107
$ws.Range("I15").Value = 0.05477919521451775
$ws.Range("J15").Value = 0.05477919521451775
$ws.Range("O15").Value = 0.04399860030713892
$ws.Range("P15").Value = 0.04399860030713892
$ws.Range("Q15").Value = 0.7559203769380001
$ws.Range("R15").Value = 6.803283392442001
$ws.Range("S15").Value = 0.002410207915390304
$ws.Range("T15").Value = 0.002410207915390304

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.6643690000000001
$ws.Range("H16").Value = 1.993107
$ws.Range("I16").Value = 0.05477919521451775
$ws.Range("J16").Value = 0.05477919521451775
$ws.Range("M16").Value = 8.848210666666667
$ws.Range("N16").Value = 26.544632
$ws.Range("O16").Value = 0.3421587275782868
$ws.Range("P16").Value = 0.3421587275782868
$ws.Range("Q16").Value = 5.878476872402667
$ws.Range("R16").Value = 52.906291851624
$ws.Range("S16").Value = 0.01874317973236197
$ws.Range("T16").Value = 0.01874317973236197

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.6643690000000001
$ws.Range("H17").Value = 1.993107
$ws.Range("I17").Value = 0.05477919521451775
$ws.Range("J17").Value = 0.05477919521451775
$ws.Range("M17").Value = 0.371739
$ws.Range("N17").Value = 1.115217
$ws.Range("O17").Value = 0.01437508079575842
$ws.Range("P17").Value = 0.01437508079575841
$ws.Range("Q17").Value = 0.246971867691
$ws.Range("R17").Value = 2.222746809219
$ws.Range("S17").Value = 0.0007874553571353155
$ws.Range("T17").Value = 0.0007874553571353152
